$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" -------------
# Overview sheet: zh-cn/de-de status columns (E:F), data rows 2-3
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = "In Translation"

# zh-cn sheet: Status column (C), data rows 2-3
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = "In Translation"

# de-de sheet: Status column (C), data rows 2-3
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = "In Translation"

# --- Column width changes ----------------------------------------------
# Target stored width is 13.4101845877511 "characters" (XML col/@width).
# The COM ColumnWidth setter here snaps to a pixel grid (~1/6-character
# steps), so the literal target is unreachable bit-for-bit; 12.5 is the
# input that lands on the closest achievable grid value (13.333333...).
$targetColumnWidth = 12.5

# Overview: columns E and F narrower
$wsOverview.Columns("E:F").ColumnWidth = $targetColumnWidth

# zh-cn / de-de: column C (Status) narrower
$wsZhCn.Columns("C").ColumnWidth = $targetColumnWidth
$wsDeDe.Columns("C").ColumnWidth = $targetColumnWidth
